# Auto-update draw results: append the 2025-10-05 Pick 3 row.
#
# The sheet stores every column as literal text (dates like "2025-10-05"
# and phase codes like "251005" are NOT real numbers/dates), so we force
# the new row to Text format before writing the values -- otherwise Excel
# would auto-convert "2025-10-05" into a date serial and "251005" into a
# number. ClearFormats() afterwards drops the now-redundant number format
# so the new cells keep the workbook's default (unstyled) appearance, just
# like every other row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 19
$rowRange = $ws.Range("A" + $newRow + ":E" + $newRow)
$rowRange.NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-10-05"
$ws.Range("B" + $newRow).Value = "Pick 3"
$ws.Range("C" + $newRow).Value = "251005"
$ws.Range("D" + $newRow).Value = "8-1-0"
$ws.Range("E" + $newRow).Value = "2025-10-05T21:34:52.320+04:00"

$rowRange.ClearFormats()
